$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param(
        [string]$MatchText,
        [string]$NewParaXml
    )

    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($MatchText)) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "NOT FOUND: $MatchText"
        return
    }

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $NewParaXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target.Range.InsertXML($pkg) | Out-Null
}

# 1) Market Basket Analysis: wrap "Apriori" with spell-check proofErr markers.
$xml1 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Market Basket Analysis for optimizing cross-selling and product recommendations in retail stores using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>Apriori</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> algorithm.</w:t></w:r></w:p>'
Replace-ParagraphXml "Market Basket Analysis" $xml1

# 2) Music Genre Classification: split the leading run into "Music " + "Genre Classification ".
$xml2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Music </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Genre Classification </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">with Support Vector Machines (SVM): </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Build a music recommendation system that classifies songs into genres based on audio features, enhancing music streaming platforms.</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:br/></w:r></w:p>'
Replace-ParagraphXml "Music Genre Classification" $xml2

# 3) LightGBM / XGBoost CTR paragraph: add strike to pPr rPr, wrap "LightGBM"/"XGBoost"/
#    "LightGBM's"/"XGBoost's" with proofErr spell-check markers.
$xml3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:u w:val="single"/></w:rPr><w:t>LightGBM</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Click-Through Rate (CTR) Prediction</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">: Utilize </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>LightGBM''s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> efficient leaf-wise tree construction and native support for categorical features to build a high-accuracy model for predicting click-through rates in online advertising, enabling better ad targeting and user experience.</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:u w:val="single"/></w:rPr><w:t>XGBoost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Click-Through Rate (CTR) Prediction</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>XGBoost''s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> ability to handle large-scale, high-dimensional data, and its capability to capture complex feature interactions make it an excellent choice for CTR prediction in online advertising.</w:t></w:r></w:p>'
Replace-ParagraphXml "LightGBM Click-Through Rate" $xml3

# 4) Facial Expression Recognition: wrap both "CapsNet" occurrences with proofErr markers.
$xml4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Facial Expression Recognition with Capsule Networks (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CapsNet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">): Apply </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CapsNet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to classify facial expressions and emotions.</w:t></w:r></w:p>'
Replace-ParagraphXml "Facial Expression Recognition" $xml4

Write-Host "Done"
